# Rename the "_old"/"_new" column-header suffixes to the respective
# input-file format-version suffixes ("_FV2404"/"_FV2410"), then turn the
# used range into a native Excel Table ("Table1") and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A..J (1-10) carry the "_old" -> "_FV2404" headers.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($baseNames[$i])_FV2404"
}

# Column K (11) is the unchanged "diff" header.

# Columns L..U (12-21) carry the "_new" -> "_FV2410" headers.
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($baseNames[$i])_FV2410"
}

# Turn the whole used range into an Excel Table named "Table1" with an
# autofilter and the (now renamed) header row as column headers.
$tableRange = $ws.Range("A1:U60")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row (split after row 1, top-left cell of the
# scrollable area is A2). Selecting A2 first makes Excel derive the
# 1-row/0-column split from the selection so the pane ends up in plain
# "frozen" state (matching Window > Freeze Panes > Freeze Top Row).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
# Restore the original active cell (A1), same as Excel leaves it after
# using the "Freeze Top Row" command from the ribbon.
$ws.Range("A1").Select() | Out-Null

Write-Output "Header rename, table creation and pane freeze applied"
